# ================================================================
# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund-holding detail, same
#    layout as the other quarterly sheets) right before "总计".
# 2) Prepend a matching summary row to the "总计" roll-up sheet.
# ================================================================

$wb = $excel.ActiveWorkbook

# --- 1) Build the new "2022-Q1" sheet ---
# Duplicate "2021-Q4" (identical column layout + styling, incl. the
# sheetPr/outline props a brand-new Worksheets.Add() sheet would
# lack) and drop it in right before "总计".
$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($totalSheet)

# NOTE: inserting a sheet "before" $totalSheet repoints that
# variable at the freshly inserted copy in this host, so every
# later use of the "总计" sheet must re-fetch it by name instead
# of reusing $totalSheet from here on.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template has 28 data rows; we only need 14, so drop the tail.
$newSheet.Rows("16:29").Delete()

# Header text (row 1 keeps the template formatting/style already)
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows. Columns B-G are text (values like "5.70" must keep
# their trailing zero, so force NumberFormat "@" before assigning,
# then clear the format back off so the cell matches a plain,
# unstyled text cell). Columns A (row index) and H (rank) are
# plain numbers.
# row 2: 011162
$newSheet.Cells.Item(2,1).Value = 0
$c = $newSheet.Cells.Item(2,2); $c.NumberFormat = "@"; $c.Value = "011162"; $c.ClearFormats()
$c = $newSheet.Cells.Item(2,3); $c.NumberFormat = "@"; $c.Value = "博时港股通领先趋势混合A"; $c.ClearFormats()
$c = $newSheet.Cells.Item(2,4); $c.NumberFormat = "@"; $c.Value = "23.56"; $c.ClearFormats()
$c = $newSheet.Cells.Item(2,5); $c.NumberFormat = "@"; $c.Value = "80.83"; $c.ClearFormats()
$c = $newSheet.Cells.Item(2,6); $c.NumberFormat = "@"; $c.Value = "5.70"; $c.ClearFormats()
$c = $newSheet.Cells.Item(2,7); $c.NumberFormat = "@"; $c.Value = "1.3429"; $c.ClearFormats()
$newSheet.Cells.Item(2,8).Value = 2
# row 3: 001605
$newSheet.Cells.Item(3,1).Value = 1
$c = $newSheet.Cells.Item(3,2); $c.NumberFormat = "@"; $c.Value = "001605"; $c.ClearFormats()
$c = $newSheet.Cells.Item(3,3); $c.NumberFormat = "@"; $c.Value = "富兰克林国海沪港深成长精选股票"; $c.ClearFormats()
$c = $newSheet.Cells.Item(3,4); $c.NumberFormat = "@"; $c.Value = "39.43"; $c.ClearFormats()
$c = $newSheet.Cells.Item(3,5); $c.NumberFormat = "@"; $c.Value = "86.45"; $c.ClearFormats()
$c = $newSheet.Cells.Item(3,6); $c.NumberFormat = "@"; $c.Value = "2.74"; $c.ClearFormats()
$c = $newSheet.Cells.Item(3,7); $c.NumberFormat = "@"; $c.Value = "1.0804"; $c.ClearFormats()
$newSheet.Cells.Item(3,8).Value = 7
# row 4: 000934
$newSheet.Cells.Item(4,1).Value = 2
$c = $newSheet.Cells.Item(4,2); $c.NumberFormat = "@"; $c.Value = "000934"; $c.ClearFormats()
$c = $newSheet.Cells.Item(4,3); $c.NumberFormat = "@"; $c.Value = "国富大中华精选混合QDII"; $c.ClearFormats()
$c = $newSheet.Cells.Item(4,4); $c.NumberFormat = "@"; $c.Value = "25.71"; $c.ClearFormats()
$c = $newSheet.Cells.Item(4,5); $c.NumberFormat = "@"; $c.Value = "83.59"; $c.ClearFormats()
$c = $newSheet.Cells.Item(4,6); $c.NumberFormat = "@"; $c.Value = "2.89"; $c.ClearFormats()
$c = $newSheet.Cells.Item(4,7); $c.NumberFormat = "@"; $c.Value = "0.7430"; $c.ClearFormats()
$newSheet.Cells.Item(4,8).Value = 7
# row 5: 006370
$newSheet.Cells.Item(5,1).Value = 3
$c = $newSheet.Cells.Item(5,2); $c.NumberFormat = "@"; $c.Value = "006370"; $c.ClearFormats()
$c = $newSheet.Cells.Item(5,3); $c.NumberFormat = "@"; $c.Value = "国富大中华精选混合QDII美元"; $c.ClearFormats()
$c = $newSheet.Cells.Item(5,4); $c.NumberFormat = "@"; $c.Value = "25.71"; $c.ClearFormats()
$c = $newSheet.Cells.Item(5,5); $c.NumberFormat = "@"; $c.Value = "83.59"; $c.ClearFormats()
$c = $newSheet.Cells.Item(5,6); $c.NumberFormat = "@"; $c.Value = "2.89"; $c.ClearFormats()
$c = $newSheet.Cells.Item(5,7); $c.NumberFormat = "@"; $c.Value = "0.7430"; $c.ClearFormats()
$newSheet.Cells.Item(5,8).Value = 7
# row 6: 007291
$newSheet.Cells.Item(6,1).Value = 4
$c = $newSheet.Cells.Item(6,2); $c.NumberFormat = "@"; $c.Value = "007291"; $c.ClearFormats()
$c = $newSheet.Cells.Item(6,3); $c.NumberFormat = "@"; $c.Value = "汇丰晋信港股通双核策略混合"; $c.ClearFormats()
$c = $newSheet.Cells.Item(6,4); $c.NumberFormat = "@"; $c.Value = "7.80"; $c.ClearFormats()
$c = $newSheet.Cells.Item(6,5); $c.NumberFormat = "@"; $c.Value = "92.66"; $c.ClearFormats()
$c = $newSheet.Cells.Item(6,6); $c.NumberFormat = "@"; $c.Value = "7.60"; $c.ClearFormats()
$c = $newSheet.Cells.Item(6,7); $c.NumberFormat = "@"; $c.Value = "0.5928"; $c.ClearFormats()
$newSheet.Cells.Item(6,8).Value = 3
# row 7: 002332
$newSheet.Cells.Item(7,1).Value = 5
$c = $newSheet.Cells.Item(7,2); $c.NumberFormat = "@"; $c.Value = "002332"; $c.ClearFormats()
$c = $newSheet.Cells.Item(7,3); $c.NumberFormat = "@"; $c.Value = "汇丰晋信沪港深股票A"; $c.ClearFormats()
$c = $newSheet.Cells.Item(7,4); $c.NumberFormat = "@"; $c.Value = "7.90"; $c.ClearFormats()
$c = $newSheet.Cells.Item(7,5); $c.NumberFormat = "@"; $c.Value = "92.60"; $c.ClearFormats()
$c = $newSheet.Cells.Item(7,6); $c.NumberFormat = "@"; $c.Value = "5.93"; $c.ClearFormats()
$c = $newSheet.Cells.Item(7,7); $c.NumberFormat = "@"; $c.Value = "0.4685"; $c.ClearFormats()
$newSheet.Cells.Item(7,8).Value = 6
# row 8: 009846
$newSheet.Cells.Item(8,1).Value = 6
$c = $newSheet.Cells.Item(8,2); $c.NumberFormat = "@"; $c.Value = "009846"; $c.ClearFormats()
$c = $newSheet.Cells.Item(8,3); $c.NumberFormat = "@"; $c.Value = "富兰克林国海港股通远见价值混合"; $c.ClearFormats()
$c = $newSheet.Cells.Item(8,4); $c.NumberFormat = "@"; $c.Value = "19.47"; $c.ClearFormats()
$c = $newSheet.Cells.Item(8,5); $c.NumberFormat = "@"; $c.Value = "86.72"; $c.ClearFormats()
$c = $newSheet.Cells.Item(8,6); $c.NumberFormat = "@"; $c.Value = "2.28"; $c.ClearFormats()
$c = $newSheet.Cells.Item(8,7); $c.NumberFormat = "@"; $c.Value = "0.4439"; $c.ClearFormats()
$newSheet.Cells.Item(8,8).Value = 9
# row 9: 011163
$newSheet.Cells.Item(9,1).Value = 7
$c = $newSheet.Cells.Item(9,2); $c.NumberFormat = "@"; $c.Value = "011163"; $c.ClearFormats()
$c = $newSheet.Cells.Item(9,3); $c.NumberFormat = "@"; $c.Value = "博时港股通领先趋势混合C"; $c.ClearFormats()
$c = $newSheet.Cells.Item(9,4); $c.NumberFormat = "@"; $c.Value = "4.68"; $c.ClearFormats()
$c = $newSheet.Cells.Item(9,5); $c.NumberFormat = "@"; $c.Value = "80.83"; $c.ClearFormats()
$c = $newSheet.Cells.Item(9,6); $c.NumberFormat = "@"; $c.Value = "5.70"; $c.ClearFormats()
$c = $newSheet.Cells.Item(9,7); $c.NumberFormat = "@"; $c.Value = "0.2668"; $c.ClearFormats()
$newSheet.Cells.Item(9,8).Value = 2
# row 10: 002333
$newSheet.Cells.Item(10,1).Value = 8
$c = $newSheet.Cells.Item(10,2); $c.NumberFormat = "@"; $c.Value = "002333"; $c.ClearFormats()
$c = $newSheet.Cells.Item(10,3); $c.NumberFormat = "@"; $c.Value = "汇丰晋信沪港深股票C"; $c.ClearFormats()
$c = $newSheet.Cells.Item(10,4); $c.NumberFormat = "@"; $c.Value = "1.23"; $c.ClearFormats()
$c = $newSheet.Cells.Item(10,5); $c.NumberFormat = "@"; $c.Value = "92.60"; $c.ClearFormats()
$c = $newSheet.Cells.Item(10,6); $c.NumberFormat = "@"; $c.Value = "5.93"; $c.ClearFormats()
$c = $newSheet.Cells.Item(10,7); $c.NumberFormat = "@"; $c.Value = "0.0729"; $c.ClearFormats()
$newSheet.Cells.Item(10,8).Value = 6
# row 11: 009017
$newSheet.Cells.Item(11,1).Value = 9
$c = $newSheet.Cells.Item(11,2); $c.NumberFormat = "@"; $c.Value = "009017"; $c.ClearFormats()
$c = $newSheet.Cells.Item(11,3); $c.NumberFormat = "@"; $c.Value = "银华港股通精选股票"; $c.ClearFormats()
$c = $newSheet.Cells.Item(11,4); $c.NumberFormat = "@"; $c.Value = "0.91"; $c.ClearFormats()
$c = $newSheet.Cells.Item(11,5); $c.NumberFormat = "@"; $c.Value = "86.12"; $c.ClearFormats()
$c = $newSheet.Cells.Item(11,6); $c.NumberFormat = "@"; $c.Value = "4.80"; $c.ClearFormats()
$c = $newSheet.Cells.Item(11,7); $c.NumberFormat = "@"; $c.Value = "0.0437"; $c.ClearFormats()
$newSheet.Cells.Item(11,8).Value = 6
# row 12: 006768
$newSheet.Cells.Item(12,1).Value = 10
$c = $newSheet.Cells.Item(12,2); $c.NumberFormat = "@"; $c.Value = "006768"; $c.ClearFormats()
$c = $newSheet.Cells.Item(12,3); $c.NumberFormat = "@"; $c.Value = "华安沪港深优选混合"; $c.ClearFormats()
$c = $newSheet.Cells.Item(12,4); $c.NumberFormat = "@"; $c.Value = "0.84"; $c.ClearFormats()
$c = $newSheet.Cells.Item(12,5); $c.NumberFormat = "@"; $c.Value = "93.09"; $c.ClearFormats()
$c = $newSheet.Cells.Item(12,6); $c.NumberFormat = "@"; $c.Value = "3.67"; $c.ClearFormats()
$c = $newSheet.Cells.Item(12,7); $c.NumberFormat = "@"; $c.Value = "0.0308"; $c.ClearFormats()
$newSheet.Cells.Item(12,8).Value = 10
# row 13: 000927
$newSheet.Cells.Item(13,1).Value = 11
$c = $newSheet.Cells.Item(13,2); $c.NumberFormat = "@"; $c.Value = "000927"; $c.ClearFormats()
$c = $newSheet.Cells.Item(13,3); $c.NumberFormat = "@"; $c.Value = "博时大中华亚太精选股票(QDII) - 美元现汇"; $c.ClearFormats()
$c = $newSheet.Cells.Item(13,4); $c.NumberFormat = "@"; $c.Value = "0.32"; $c.ClearFormats()
$c = $newSheet.Cells.Item(13,5); $c.NumberFormat = "@"; $c.Value = "92.94"; $c.ClearFormats()
$c = $newSheet.Cells.Item(13,6); $c.NumberFormat = "@"; $c.Value = "4.03"; $c.ClearFormats()
$c = $newSheet.Cells.Item(13,7); $c.NumberFormat = "@"; $c.Value = "0.0129"; $c.ClearFormats()
$newSheet.Cells.Item(13,8).Value = 9
# row 14: 050015
$newSheet.Cells.Item(14,1).Value = 12
$c = $newSheet.Cells.Item(14,2); $c.NumberFormat = "@"; $c.Value = "050015"; $c.ClearFormats()
$c = $newSheet.Cells.Item(14,3); $c.NumberFormat = "@"; $c.Value = "博时大中华亚太精选股票(QDII) -人民币"; $c.ClearFormats()
$c = $newSheet.Cells.Item(14,4); $c.NumberFormat = "@"; $c.Value = "0.32"; $c.ClearFormats()
$c = $newSheet.Cells.Item(14,5); $c.NumberFormat = "@"; $c.Value = "92.94"; $c.ClearFormats()
$c = $newSheet.Cells.Item(14,6); $c.NumberFormat = "@"; $c.Value = "4.03"; $c.ClearFormats()
$c = $newSheet.Cells.Item(14,7); $c.NumberFormat = "@"; $c.Value = "0.0129"; $c.ClearFormats()
$newSheet.Cells.Item(14,8).Value = 9
# row 15: 001824
$newSheet.Cells.Item(15,1).Value = 13
$c = $newSheet.Cells.Item(15,2); $c.NumberFormat = "@"; $c.Value = "001824"; $c.ClearFormats()
$c = $newSheet.Cells.Item(15,3); $c.NumberFormat = "@"; $c.Value = "博时沪港深成长企业混合"; $c.ClearFormats()
$c = $newSheet.Cells.Item(15,4); $c.NumberFormat = "@"; $c.Value = "0.15"; $c.ClearFormats()
$c = $newSheet.Cells.Item(15,5); $c.NumberFormat = "@"; $c.Value = "93.41"; $c.ClearFormats()
$c = $newSheet.Cells.Item(15,6); $c.NumberFormat = "@"; $c.Value = "6.52"; $c.ClearFormats()
$c = $newSheet.Cells.Item(15,7); $c.NumberFormat = "@"; $c.Value = "0.0098"; $c.ClearFormats()
$newSheet.Cells.Item(15,8).Value = 1

# --- 2) Prepend a "2022-Q1" row to the "总计" roll-up sheet ---
# Re-fetch by name (see note above) to get the real "总计" sheet.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Inherit the index-column style (s="2") from the row pushed down.
$totalSheet.Cells.Item(3,1).Copy()
$totalSheet.Cells.Item(2,1).PasteSpecial(-4122)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 14
$totalSheet.Cells.Item(2,4).Value = 5.86
$totalSheet.Range("B2:D2").ClearFormats()

# Renumber the index column (A) for the rows pushed down by the
# insert, so it stays a contiguous 0..n sequence.
$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(5,1).Value = 3
$totalSheet.Cells.Item(6,1).Value = 4
$totalSheet.Cells.Item(7,1).Value = 5

